# Update crypto price/volume data per the Sep 26 2023 GitHub Actions scrape refresh.
# Row order for a couple of coin pairs changed (rank swap), and prices/volumes were refreshed.
# A leading apostrophe is used for Price values that would otherwise be auto-parsed as numbers
# by Excel, so they stay stored as text (matching the sheet's existing text-based Price column).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "26.372.80"
$ws.Range('E2').Value = "  +0.41%  "
$ws.Range('D3').Value = "1.592.59"
$ws.Range('E3').Value = "  +0.66%  "
$ws.Range('E4').Value = "  -0.31%  "
$ws.Range('D5').Value = "'211.45"
$ws.Range('E5').Value = "  +0.89%  "
$ws.Range('E6').Value = "  +0.43%  "
$ws.Range('E7').Value = "  -0.30%  "
$ws.Range('E8').Value = "  +0.60%  "
$ws.Range('E9').Value = "  -0.01%  "
$ws.Range('D10').Value = "'19.49"
$ws.Range('E10').Value = "  -0.50%  "
$ws.Range('D11').Value = "'0.0847"
$ws.Range('E11').Value = "  +0.23%  "
$ws.Range('D12').Value = "1.816.09"
$ws.Range('E12').Value = "  +0.63%  "
$ws.Range('B13').Value = "Polkadot"
$ws.Range('C13').Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range('D13').Value = "'4.06"
$ws.Range('E13').Value = "  +1.08%  "
$ws.Range('B14').Value = "WrappedEther"
$ws.Range('C14').Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D14').Value = "1.579.82"
$ws.Range('E14').Value = "  -0.07%  "
$ws.Range('E15').Value = "  +1.42%  "
$ws.Range('D16').Value = "'64.87"
$ws.Range('E16').Value = "  +0.50%  "
$ws.Range('D17').Value = "26.370.31"
$ws.Range('E18').Value = "  -0.96%  "
$ws.Range('D20').Value = "'212.26"
$ws.Range('E20').Value = "  +2.76%  "
$ws.Range('E21').Value = "  -0.31%  "
$ws.Range('E22').Value = "  +1.27%  "
$ws.Range('B23').Value = "Toncoin"
$ws.Range('C23').Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range('D23').Value = "'2.18"
$ws.Range('E23').Value = "  -1.75%  "
$ws.Range('B24').Value = "Avalanche"
$ws.Range('C24').Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range('D24').Value = "'9.03"
$ws.Range('E24').Value = "  +2.29%  "
$ws.Range('D25').Value = "'144.59"
$ws.Range('E25').Value = "  +0.00%  "
$ws.Range('E26').Value = "  -0.25%  "
$ws.Range('E27').Value = "  +0.88%  "
$ws.Range('E28').Value = "  -0.49%  "
$ws.Range('D29').Value = "'15.24"
$ws.Range('E29').Value = "  -0.02%  "
$ws.Range('E30').Value = "  +0.12%  "
$ws.Range('E31').Value = "  +1.09%  "
$ws.Range('E32').Value = "  -0.22%  "
$ws.Range('E33').Value = "  +1.13%  "
$ws.Range('D34').Value = "1.338.65"
$ws.Range('E34').Value = "  +4.28%  "
$ws.Range('E35').Value = "  -1.33%  "
$ws.Range('E36').Value = "  -0.54%  "
$ws.Range('D37').Value = "'1.48"
$ws.Range('E37').Value = "  +0.06%  "
$ws.Range('E38').Value = "  +0.04%  "
$ws.Range('B39').Value = "WEMIXToken"
$ws.Range('C39').Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range('D39').Value = "'1.06"
$ws.Range('E39').Value = "  -16.01%  "
$ws.Range('B40').Value = "ARBITRUM"
$ws.Range('C40').Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range('D40').Value = "'0.819"
$ws.Range('E40').Value = "  +0.44%  "
$ws.Range('D41').Value = "'5.77"
$ws.Range('E41').Value = "  +5.14%  "
$ws.Range('E42').Value = "  -0.30%  "
$ws.Range('E43').Value = "  +0.31%  "
$ws.Range('D44').Value = "'0.764"
$ws.Range('E44').Value = "  -0.71%  "
$ws.Range('D45').Value = "1.728.76"
$ws.Range('E45').Value = "  +0.62%  "
$ws.Range('D46').Value = "'61.92"
$ws.Range('E46').Value = "  -0.61%  "
$ws.Range('E47').Value = "  -0.73%  "
$ws.Range('E48').Value = "  -3.40%  "
$ws.Range('D49').Value = "'0.0986"
$ws.Range('E49').Value = "  -3.08%  "
$ws.Range('D50').Value = "'0.0506"
$ws.Range('E50').Value = "  -0.63%  "
$ws.Range('E51').Value = "  -0.32%  "
